$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17,8).Value = 1000380.75
$ws.Cells.Item(17,10).Value = 1000380.75
$ws.Cells.Item(17,12).Value = 3001142.25
$ws.Cells.Item(17,14).Value = -3001478.25
$ws.Cells.Item(86,8).Value = 80999.2
$ws.Cells.Item(86,9).Value = 1496
$ws.Cells.Item(86,11).Value = 1496
$ws.Cells.Item(86,13).Value = -373
$ws.Cells.Item(89,8).Value = 80999.2
$ws.Cells.Item(89,9).Value = 1496
$ws.Cells.Item(89,11).Value = 7480
$ws.Cells.Item(89,13).Value = -1864
$ws.Cells.Item(112,8).Value = 1161.0571
$ws.Cells.Item(112,10).Value = 1160.5312
$ws.Cells.Item(112,12).Value = 3481.5936
$ws.Cells.Item(112,14).Value = -5697.5936
$ws.Cells.Item(137,8).Value = 1052.1765
$ws.Cells.Item(137,9).Value = 993.0606
$ws.Cells.Item(137,11).Value = 2979.1818
$ws.Cells.Item(137,13).Value = -429.1818000000003
$ws.Cells.Item(138,8).Value = 2624.307
$ws.Cells.Item(138,9).Value = 1274.7736
$ws.Cells.Item(138,10).Value = 4667.8857
$ws.Cells.Item(138,11).Value = 3824.3208
$ws.Cells.Item(138,12).Value = 14003.6571
$ws.Cells.Item(138,13).Value = 1315.6792
$ws.Cells.Item(138,14).Value = -24283.6571
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32,8).Value = 1295.31
$ws.Cells.Item(32,9).Value = 1262.4183
$ws.Cells.Item(32,10).Value = 2907
$ws.Cells.Item(32,11).Value = 1262.4183
$ws.Cells.Item(32,12).Value = 2907
$ws.Cells.Item(32,13).Value = -975.4183
$ws.Cells.Item(32,14).Value = -3481
$ws.Cells.Item(45,8).Value = 1574.6875
$ws.Cells.Item(45,9).Value = 1245.1111
$ws.Cells.Item(45,11).Value = 1245.1111
$ws.Cells.Item(45,13).Value = -868.1111000000001
$ws.Cells.Item(61,8).Value = 546.64105
$ws.Cells.Item(61,9).Value = 530.2432
$ws.Cells.Item(61,10).Value = 850
$ws.Cells.Item(61,11).Value = 530.2432
$ws.Cells.Item(61,12).Value = 850
$ws.Cells.Item(61,13).Value = -318.2432
$ws.Cells.Item(61,14).Value = -1274
$ws.Cells.Item(74,8).Value = 2711.034
$ws.Cells.Item(74,9).Value = 3078.88
$ws.Cells.Item(74,10).Value = 667.44446
$ws.Cells.Item(74,11).Value = 3078.88
$ws.Cells.Item(74,12).Value = 667.44446
$ws.Cells.Item(74,13).Value = -2204.88
$ws.Cells.Item(74,14).Value = -2415.44446
$ws.Cells.Item(77,8).Value = 2711.034
$ws.Cells.Item(77,9).Value = 3078.88
$ws.Cells.Item(77,10).Value = 667.44446
$ws.Cells.Item(77,11).Value = 15394.4
$ws.Cells.Item(77,12).Value = 3337.2223
$ws.Cells.Item(77,13).Value = -11026.4
$ws.Cells.Item(77,14).Value = -12073.2223
$ws.Cells.Item(122,8).Value = 1734.3793
$ws.Cells.Item(122,9).Value = 1396.0385
$ws.Cells.Item(122,10).Value = 4666.6665
$ws.Cells.Item(122,11).Value = 4188.1155
$ws.Cells.Item(122,12).Value = 13999.9995
$ws.Cells.Item(122,13).Value = -1738.1155
$ws.Cells.Item(122,14).Value = -18899.9995
$ws.Cells.Item(132,8).Value = 1483.1466
$ws.Cells.Item(132,9).Value = 747.2857
$ws.Cells.Item(132,10).Value = 3652
$ws.Cells.Item(132,11).Value = 2241.8571
$ws.Cells.Item(132,12).Value = 10956
$ws.Cells.Item(132,13).Value = 288.1428999999998
$ws.Cells.Item(132,14).Value = -16016
$ws.Cells.Item(136,8).Value = 546.64105
$ws.Cells.Item(136,9).Value = 530.2432
$ws.Cells.Item(136,10).Value = 850
$ws.Cells.Item(136,11).Value = 1590.7296
$ws.Cells.Item(136,12).Value = 2550
$ws.Cells.Item(136,13).Value = 959.2703999999999
$ws.Cells.Item(136,14).Value = -7650
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20,8).Value = 27229.5
$ws.Cells.Item(20,9).Value = 0
$ws.Cells.Item(20,10).Value = 27229.5
$ws.Cells.Item(20,11).Value = 0
$ws.Cells.Item(20,12).ClearContents()
$ws.Cells.Item(20,13).Value = 27229.5
$ws.Cells.Item(20,14).Value = -27723.5
$ws.Cells.Item(134,8).Value = 1183.807
$ws.Cells.Item(134,9).Value = 679.6739
$ws.Cells.Item(134,10).Value = 3292
$ws.Cells.Item(134,11).Value = 2039.0217
$ws.Cells.Item(134,12).Value = 9876
$ws.Cells.Item(134,13).Value = 495.9783
$ws.Cells.Item(134,14).Value = -14946
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4,8).Value = 3656667.2
$ws.Cells.Item(4,9).Value = 10000000
$ws.Cells.Item(4,10).Value = 485001
$ws.Cells.Item(4,11).Value = 10000000
$ws.Cells.Item(4,12).Value = 485001
$ws.Cells.Item(4,13).Value = -9999888
$ws.Cells.Item(4,14).Value = -485225
$ws.Cells.Item(16,8).Value = 3328.1
$ws.Cells.Item(16,9).Value = 2586.7778
$ws.Cells.Item(16,10).Value = 10000
$ws.Cells.Item(16,11).Value = 2586.7778
$ws.Cells.Item(16,12).Value = 10000
$ws.Cells.Item(16,13).Value = -2299.7778
$ws.Cells.Item(16,14).Value = -10574
$ws.Cells.Item(95,8).Value = 25000
$ws.Cells.Item(95,10).Value = 25000
$ws.Cells.Item(95,12).Value = 25000
$ws.Cells.Item(95,14).Value = -30492
$ws.Cells.Item(113,8).Value = 3328.1
$ws.Cells.Item(113,9).Value = 2586.7778
$ws.Cells.Item(113,10).Value = 10000
$ws.Cells.Item(113,11).Value = 2586.7778
$ws.Cells.Item(113,12).Value = 10000
$ws.Cells.Item(113,13).Value = -416.7777999999998
$ws.Cells.Item(113,14).Value = -14340
$ws.Cells.Item(134,8).Value = 1046.5294
$ws.Cells.Item(134,9).Value = 1046.5294
$ws.Cells.Item(134,10).Value = 0
$ws.Cells.Item(134,11).Value = 3139.5882
$ws.Cells.Item(134,12).Value = 0
$ws.Cells.Item(134,13).ClearContents()
$ws.Cells.Item(134,14).Value = -604.5881999999997
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5,8).Value = 549322.9399999999
$ws.Cells.Item(5,9).Value = 716.75
$ws.Cells.Item(5,10).Value = 1097929.1
$ws.Cells.Item(5,11).Value = 2150.25
$ws.Cells.Item(5,12).Value = 3293787.3
$ws.Cells.Item(5,13).Value = -2038.25
$ws.Cells.Item(5,14).Value = -3294011.3
$ws.Cells.Item(122,8).Value = 871.2
$ws.Cells.Item(122,9).Value = 589.25
$ws.Cells.Item(122,10).Value = 1999
$ws.Cells.Item(122,11).Value = 5303.25
$ws.Cells.Item(122,12).Value = 17991
$ws.Cells.Item(122,13).Value = -2853.25
$ws.Cells.Item(122,14).Value = -22891
$ws.Cells.Item(132,8).Value = 1887.25
$ws.Cells.Item(132,9).Value = 2717.6667
$ws.Cells.Item(132,10).Value = 1531.3572
$ws.Cells.Item(132,11).Value = 24459.0003
$ws.Cells.Item(132,12).Value = 13782.2148
$ws.Cells.Item(132,13).Value = -21929.0003
$ws.Cells.Item(132,14).Value = -18842.2148
$ws.Cells.Item(135,8).Value = 549322.9399999999
$ws.Cells.Item(135,9).Value = 716.75
$ws.Cells.Item(135,10).Value = 1097929.1
$ws.Cells.Item(135,11).Value = 6450.75
$ws.Cells.Item(135,12).Value = 9881361.9
$ws.Cells.Item(135,13).Value = -3915.75
$ws.Cells.Item(135,14).Value = -9886431.9
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70,8).Value = 6032.68
$ws.Cells.Item(70,9).Value = 4822.9414
$ws.Cells.Item(70,10).Value = 8603.375
$ws.Cells.Item(70,11).Value = 4822.9414
$ws.Cells.Item(70,12).Value = 8603.375
$ws.Cells.Item(70,13).Value = -4552.9414
$ws.Cells.Item(70,14).Value = -9143.375
$ws.Cells.Item(73,8).Value = 6032.68
$ws.Cells.Item(73,9).Value = 4822.9414
$ws.Cells.Item(73,10).Value = 8603.375
$ws.Cells.Item(73,11).Value = 4822.9414
$ws.Cells.Item(73,12).Value = 8603.375
$ws.Cells.Item(73,13).Value = -3886.9414
$ws.Cells.Item(73,14).Value = -10475.375
$ws.Cells.Item(80,8).Value = 3102.5
$ws.Cells.Item(80,9).Value = 2227.5
$ws.Cells.Item(80,10).Value = 3977.5
$ws.Cells.Item(80,11).Value = 2227.5
$ws.Cells.Item(80,12).Value = 3977.5
$ws.Cells.Item(80,13).Value = -1229.5
$ws.Cells.Item(80,14).Value = -5973.5
$ws.Cells.Item(83,8).Value = 3102.5
$ws.Cells.Item(83,9).Value = 2227.5
$ws.Cells.Item(83,10).Value = 3977.5
$ws.Cells.Item(83,11).Value = 11137.5
$ws.Cells.Item(83,12).Value = 19887.5
$ws.Cells.Item(83,13).Value = -6145.5
$ws.Cells.Item(83,14).Value = -29871.5
$ws.Cells.Item(93,8).Value = 20100
$ws.Cells.Item(93,10).Value = 20100
$ws.Cells.Item(93,12).Value = 20100
$ws.Cells.Item(93,14).Value = -23844
$ws.Cells.Item(97,8).Value = 1363.3846
$ws.Cells.Item(97,9).Value = 1523.1818
$ws.Cells.Item(97,10).Value = 484.5
$ws.Cells.Item(97,11).Value = 1523.1818
$ws.Cells.Item(97,12).Value = 484.5
$ws.Cells.Item(97,13).Value = -1027.1818
$ws.Cells.Item(97,14).Value = -1476.5
$ws.Cells.Item(122,8).Value = 2330.158
$ws.Cells.Item(122,9).Value = 2002.64
$ws.Cells.Item(122,10).Value = 2960
$ws.Cells.Item(122,11).Value = 6007.92
$ws.Cells.Item(122,12).Value = 8880
$ws.Cells.Item(122,13).Value = -3557.92
$ws.Cells.Item(122,14).Value = -13780
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46,8).Value = 168800.5
$ws.Cells.Item(46,9).Value = 202460.2
$ws.Cells.Item(46,10).Value = 502
$ws.Cells.Item(46,11).Value = 202460.2
$ws.Cells.Item(46,12).Value = 502
$ws.Cells.Item(46,13).Value = -202272.2
$ws.Cells.Item(46,14).Value = -878
$ws.Cells.Item(55,8).Value = 1199.25
$ws.Cells.Item(55,9).Value = 844.55554
$ws.Cells.Item(55,10).Value = 2263.3333
$ws.Cells.Item(55,11).Value = 844.55554
$ws.Cells.Item(55,12).Value = 2263.3333
$ws.Cells.Item(55,13).Value = -671.55554
$ws.Cells.Item(55,14).Value = -2609.3333
$ws.Cells.Item(61,8).Value = 1753
$ws.Cells.Item(61,9).Value = 1669
$ws.Cells.Item(61,11).Value = 1669
$ws.Cells.Item(61,13).Value = -1467
$ws.Cells.Item(82,8).Value = 997.9167
$ws.Cells.Item(82,9).Value = 922.3125
$ws.Cells.Item(82,10).Value = 1149.125
$ws.Cells.Item(82,11).Value = 922.3125
$ws.Cells.Item(82,12).Value = 1149.125
$ws.Cells.Item(82,13).Value = -561.3125
$ws.Cells.Item(82,14).Value = -1871.125
$ws.Cells.Item(85,8).Value = 997.9167
$ws.Cells.Item(85,9).Value = 922.3125
$ws.Cells.Item(85,10).Value = 1149.125
$ws.Cells.Item(85,11).Value = 922.3125
$ws.Cells.Item(85,12).Value = 1149.125
$ws.Cells.Item(85,13).Value = 325.6875
$ws.Cells.Item(85,14).Value = -3645.125
$ws.Cells.Item(113,8).Value = 1753
$ws.Cells.Item(113,9).Value = 1669
$ws.Cells.Item(113,11).Value = 1669
$ws.Cells.Item(113,13).Value = 501
$ws.Cells.Item(122,8).Value = 2095.5
$ws.Cells.Item(122,9).Value = 1402
$ws.Cells.Item(122,10).Value = 2326.6667
$ws.Cells.Item(122,11).Value = 4206
$ws.Cells.Item(122,12).Value = 6980.000100000001
$ws.Cells.Item(122,13).Value = -1756
$ws.Cells.Item(122,14).Value = -11880.0001
$ws.Cells.Item(132,8).Value = 5476.1123
$ws.Cells.Item(132,9).Value = 5937.521
$ws.Cells.Item(132,10).Value = 4784
$ws.Cells.Item(132,11).Value = 17812.563
$ws.Cells.Item(132,12).Value = 14352
$ws.Cells.Item(132,13).Value = -15282.563
$ws.Cells.Item(132,14).Value = -19412
$ws.Cells.Item(136,8).Value = 2093.3333
$ws.Cells.Item(136,9).Value = 2374.375
$ws.Cells.Item(136,10).Value = 969.1667
$ws.Cells.Item(136,11).Value = 7123.125
$ws.Cells.Item(136,12).Value = 2907.5001
$ws.Cells.Item(136,13).Value = -4573.125
$ws.Cells.Item(136,14).Value = -8007.5001
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122,8).Value = 51854.4
$ws.Cells.Item(122,9).Value = 68232.53
$ws.Cells.Item(122,10).Value = 2720
$ws.Cells.Item(122,11).Value = 204697.59
$ws.Cells.Item(122,12).Value = 8160
$ws.Cells.Item(122,13).Value = -202247.59
$ws.Cells.Item(122,14).Value = -13060
$ws.Cells.Item(133,8).Value = 30000
$ws.Cells.Item(133,10).Value = 30000
$ws.Cells.Item(133,12).Value = 30000
$ws.Cells.Item(133,14).Value = -40120
$ws.Cells.Item(136,8).Value = 1070.4
$ws.Cells.Item(136,9).Value = 555.7805
$ws.Cells.Item(136,10).Value = 2577.5
$ws.Cells.Item(136,11).Value = 1667.3415
$ws.Cells.Item(136,12).Value = 7732.5
$ws.Cells.Item(136,13).Value = 882.6585
$ws.Cells.Item(136,14).Value = -12832.5
